$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "CustomerData"

# Header row (bold)
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "PostCode"
$ws.Range("A1:C1").Font.Bold = $true

# Data rows
$ws.Range("A2").Value = "Meenakshi"
$ws.Range("B2").Value = "Rana"
$ws.Range("C2").Value = 111

$ws.Range("A3").Value = "Vamika "
$ws.Range("B3").Value = "Dogra"
$ws.Range("C3").Value = 111

$ws.Range("A4").Value = "Pardeep "
$ws.Range("B4").Value = "Dogra"
$ws.Range("C4").Value = 222

# Selection matches the authored sheet view
$ws.Range("A1:C1").Select()

# Page setup (portrait) matches the authored sheet
$ws.PageSetup.Orientation = 1
